$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Aké"
$ws.Range("B2").Value = 81
$ws.Range("E2").Value = "LB"
$ws.Range("G2").Value = "Netherlands"
$ws.Range("H2").Value = "Premier League"
$ws.Range("I2").Value = "Manchester City"
$ws.Range("N2").Value = 700
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 209

# Row 3
$ws.Range("A3").Value = "Danjuma"
$ws.Range("B3").Value = 81
$ws.Range("D3").Value = "Rare"
$ws.Range("E3").Value = "LM"
$ws.Range("G3").Value = "Netherlands"
$ws.Range("H3").Value = "Premier League"
$ws.Range("I3").Value = "Everton"
$ws.Range("P3").Value = 239

# Row 4
$ws.Range("A4").Value = "Dalot Teixeira"
$ws.Range("B4").Value = 80
$ws.Range("E4").Value = "RB"
$ws.Range("G4").Value = "Portugal"
$ws.Range("H4").Value = "Premier League"
$ws.Range("I4").Value = "Manchester Utd"
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 284

# Row 5
$ws.Range("A5").Value = "García Córdoba"
$ws.Range("B5").Value = 80
$ws.Range("E5").Value = "ST"
$ws.Range("G5").Value = "Spain"
$ws.Range("H5").Value = "Barclays WSL"
$ws.Range("I5").Value = "Manchester Utd"
$ws.Range("N5").Value = 750
$ws.Range("P5").Value = 297

# Row 6
$ws.Range("A6").Value = "Elustondo"
$ws.Range("B6").Value = 80
$ws.Range("D6").Value = "Common"
$ws.Range("E6").Value = "CB"
$ws.Range("G6").Value = "Spain"
$ws.Range("H6").Value = "LALIGA EA SPORTS"
$ws.Range("I6").Value = "Real Sociedad"
$ws.Range("N6").Value = 550
$ws.Range("O6").Value = 3
$ws.Range("P6").Value = 303

# Row 7
$ws.Range("A7").Value = "Martial"
$ws.Range("E7").Value = "ST"
$ws.Range("I7").Value = "Manchester Utd"
$ws.Range("P7").Value = 310

# Row 8
$ws.Range("A8").Value = "Catena Marugán"
$ws.Range("B8").Value = 79
$ws.Range("D8").Value = "Common"
$ws.Range("E8").Value = "CB"
$ws.Range("G8").Value = "Spain"
$ws.Range("I8").Value = "CA Osasuna"
$ws.Range("N8").Value = 500
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 332

# Row 9
$ws.Range("A9").Value = "Pinillos Moreno"
$ws.Range("E9").Value = "RM"
$ws.Range("G9").Value = "Spain"
$ws.Range("H9").Value = "Liga F"
$ws.Range("I9").Value = "Madrid CFF"
$ws.Range("N9").Value = 500
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 334

# Row 10
$ws.Range("A10").Value = "McTominay"
$ws.Range("B10").Value = 79
$ws.Range("E10").Value = "CM"
$ws.Range("G10").Value = "Scotland"
$ws.Range("H10").Value = "Premier League"
$ws.Range("I10").Value = "Manchester Utd"
$ws.Range("N10").Value = 500
$ws.Range("P10").Value = 336

# Row 11
$ws.Range("A11").Value = "Herrera Pirón"
$ws.Range("B11").Value = 79
$ws.Range("E11").Value = "GK"
$ws.Range("G11").Value = "Spain"
$ws.Range("H11").Value = "LALIGA EA SPORTS"
$ws.Range("I11").Value = "CA Osasuna"
$ws.Range("N11").Value = 500
$ws.Range("P11").Value = 351

# Row 12
$ws.Range("A12").Value = "Maitane"
$ws.Range("B12").Value = 79
$ws.Range("C12").Value = "Gold"
$ws.Range("E12").Value = "CM"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = "Spain"
$ws.Range("H12").Value = "NWSL"
$ws.Range("I12").Value = "NJ/NY Gotham"
$ws.Range("N12").Value = 500
$ws.Range("O12").Value = 2
$ws.Range("P12").Value = 374
